{"js": "// Helper to find the first search hit for a unique string and replace its text.\nasync function replaceOnce(context, searchText, newText) {\n  const results = context.document.body.search(searchText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + searchText);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 1) Salutation line: \"Dear Mr./Mrs./Ms. [Last Name],\" -> \"Dear Capstone Engineering Society Committee,\"\nawait replaceOnce(\n  context,\n  \"Dear Mr./Mrs./Ms. [Last Name],\",\n  \"Dear Capstone Engineering Society Committee,\"\n);\n\n// 2) Opening paragraph: purpose of the letter changes from \"admission\" to the award nomination.\nawait replaceOnce(\n  context,\n  \"for admission to your graduate program.\",\n  \"for the Capstone Engineering Society Outstanding Senior Award.\"\n);\n\n// 3) Spelling fix: \"labelling\" -> \"labeling\"\nawait replaceOnce(context, \"labelling\", \"labeling\");\n\n// 4) \"the high-performance computer\" -> \"a high-performance computer\"\nawait replaceOnce(\n  context,\n  \"how to use the high-performance computer.\",\n  \"how to use a high-performance computer.\"\n);\n\n// 5) Insert a brand-new paragraph describing Kate's continued use of automation skills,\n//    right before the concluding paragraph, and rewrite that concluding paragraph.\nconst oldConclusion =\n  \"In conclusion, I confidently recommend Kate Sanborn as a student in your graduate program. Based on the work she has completed this summer and the growth that she has had, I know she will be an excellent addition to your student body.\";\n\nconst concludingResults = context.document.body.search(oldConclusion, { matchCase: true, matchWholeWord: false });\nconcludingResults.load(\"items\");\nawait context.sync();\nif (concludingResults.items.length === 0) {\n  throw new Error(\"Concluding paragraph not found\");\n}\nconst concludingRange = concludingResults.items[0];\n\nconst newMiddleParagraphText =\n  \"Since the CAT Vehicle REU program, Kate has continued to use what she learned about automation in other work experiences and the classroom. Last summer, she worked as an IT Process Automation intern at Ulta Beauty. In this role, she used Blue Prism to automate tedious manual processes, saving time and resources. Currently, her senior design team is building a robot for the IEEE SoutheastCon 2023 conference hardware competition. Using what she learned from the REU, Kate is writing code for path planning and mission control to allow the robot to run autonomously. This spring, she also joined the Crimson Autonomous Kart team through the Autonomous Vehicles class, a computer science elective. She continues to use what she learned from research in her role on this team.\";\n\n// NOTE: insertParagraph(before) re-targets the original range object to the\n// newly inserted paragraph, so we must re-search for the conclusion text\n// afterwards to get a fresh reference to it before editing it further.\nconcludingRange.insertParagraph(newMiddleParagraphText, Word.InsertLocation.before);\nawait context.sync();\n\nconst newConclusionText =\n  \"In conclusion, I confidently recommend Kate Sanborn for the Capstone Engineering Society Outstanding Senior Award. Based on the work she completed that summer and the growth that she had, I know she will be an excellent candidate for this award.\";\n\nconst concludingResults2 = context.document.body.search(oldConclusion, { matchCase: true, matchWholeWord: false });\nconcludingResults2.load(\"items\");\nawait context.sync();\nif (concludingResults2.items.length === 0) {\n  throw new Error(\"Concluding paragraph not found after insertion\");\n}\nconcludingResults2.items[0].insertText(newConclusionText, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Word COM interop script to apply the \"Capstone Engineering Society Outstanding\n# Senior Award\" rewrite to the recommendation letter.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($searchText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $searchText\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = 0\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.Execute($find.Text, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n\nfunction Find-Range($searchText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $searchText\n    $find.Forward = $true\n    $find.Wrap = 0\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $found = $find.Execute()\n    if (-not $found) {\n        throw \"Text not found: $searchText\"\n    }\n    return $find.Parent\n}\n\n# 1) Salutation line.\nReplace-Text \"Dear Mr./Mrs./Ms. [Last Name],\" \"Dear Capstone Engineering Society Committee,\"\n\n# 2) Opening paragraph: purpose of the letter changes from \"admission\" to the award nomination.\nReplace-Text \"for admission to your graduate program.\" \"for the Capstone Engineering Society Outstanding Senior Award.\"\n\n# 3) Spelling fix.\nReplace-Text \"labelling\" \"labeling\"\n\n# 4) \"the high-performance computer\" -> \"a high-performance computer\"\nReplace-Text \"how to use the high-performance computer.\" \"how to use a high-performance computer.\"\n\n# 5) Insert a brand-new paragraph describing Kate's continued use of automation skills,\n#    right before the concluding paragraph, then rewrite that concluding paragraph.\n$oldConclusion = \"In conclusion, I confidently recommend Kate Sanborn as a student in your graduate program. Based on the work she has completed this summer and the growth that she has had, I know she will be an excellent addition to your student body.\"\n\n$concludingRange = Find-Range $oldConclusion\n$concludingRange.InsertParagraphBefore()\n\n# InsertParagraphBefore() can leave stale Start/End offsets cached on the\n# range/find objects involved, so re-find the conclusion text to get a\n# reliable, up-to-date range before doing any further positional math.\n$freshConcludingRange = Find-Range $oldConclusion\n\n# The newly inserted blank paragraph sits immediately before the (now\n# shifted) conclusion paragraph; grab it via a collapsed range just before\n# the conclusion text starts.\n$blankRange = $d.Range($freshConcludingRange.Start - 1, $freshConcludingRange.Start - 1)\n$blankPara = $blankRange.Paragraphs(1)\n$newMiddleParagraphText = \"Since the CAT Vehicle REU program, Kate has continued to use what she learned about automation in other work experiences and the classroom. Last summer, she worked as an IT Process Automation intern at Ulta Beauty. In this role, she used Blue Prism to automate tedious manual processes, saving time and resources. Currently, her senior design team is building a robot for the IEEE SoutheastCon 2023 conference hardware competition. Using what she learned from the REU, Kate is writing code for path planning and mission control to allow the robot to run autonomously. This spring, she also joined the Crimson Autonomous Kart team through the Autonomous Vehicles class, a computer science elective. She continues to use what she learned from research in her role on this team.\"\n$blankPara.Range.Text = $newMiddleParagraphText\n\n# Filling in the blank paragraph's text shifts offsets again, so re-find the\n# conclusion text one more time before rewriting it.\n$finalConcludingRange = Find-Range $oldConclusion\n$newConclusionText = \"In conclusion, I confidently recommend Kate Sanborn for the Capstone Engineering Society Outstanding Senior Award. Based on the work she completed that summer and the growth that she had, I know she will be an excellent candidate for this award.\"\n$finalConcludingRange.Text = $newConclusionText\n"}
